# Daily attendance processing - 2025-10-28 17:46:39
# Normalizes the "Recorded By" column (G): the audit-log order for the
# recorded-by list is reversed (most-recent-first -> most-recent-last),
# except for a handful of dates whose entries were already reconciled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skipDates = @("27/09/2025", "08/10/2025", "11/10/2025", "14/10/2025", "16/10/2025")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dateVal = $ws.Cells.Item($r, 5).Value2
    $recordedBy = $ws.Cells.Item($r, 7).Value2

    if ([string]::IsNullOrEmpty($recordedBy)) { continue }
    if (-not $recordedBy.StartsWith("System, ")) { continue }
    if ($skipDates -contains $dateVal) { continue }

    $parts = $recordedBy.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $reversed = @()
    for ($i = $parts.Length - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $ws.Cells.Item($r, 7).Value2 = [string]::Join(", ", $reversed)
}
